# Updated README. Changed excel importing to take relative references.
# Replace absolute Windows file paths in the "atoms" column (J) with
# relative paths for the H2 and H2O CONTCAR references.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = ".\H2\CONTCAR"
$ws.Range("J4").Value = ".\H2O\CONTCAR"

# Leave the final selection on F2, matching the saved file's cursor position.
$ws.Range("F2").Select()
